$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3"
#    worksheet (so it becomes the 2nd tab, right after "总计").
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch a handle to the (still existing) "2022-Q3" sheet by name - its
# index shifted by one once the new sheet was inserted in front of it.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Copy the header row (with its styling) from the 2022-Q3 sheet so the new
# sheet matches the look of all the other quarterly sheets.
$q3Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Copy the styling used for the numeric row-index column (column A) as well.
$q3Sheet.Range("A2").Copy($newSheet.Range("A2:A10"))

# ---------------------------------------------------------------------
# Fill in the 2022-Q4 fund-holding data.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$rows = @(
    @{ idx = 0; code = "161724"; name = "招商中证煤炭等权指数（LOF）A";        size = "17.24"; pos = "93.84"; pct = "3.12"; mv = "0.5379"; rank = 2 },
    @{ idx = 1; code = "501029"; name = "华宝标普中国A股红利机会指数（LOF）A"; size = "10.97"; pos = "94.25"; pct = "1.86"; mv = "0.2040"; rank = 3 },
    @{ idx = 2; code = "005125"; name = "华宝标普中国A股红利机会指数C";        size = "3.29";  pos = "94.25"; pct = "1.86"; mv = "0.0612"; rank = 3 },
    @{ idx = 3; code = "013596"; name = "招商中证煤炭等权指数（LOF）C";        size = "1.56";  pos = "93.84"; pct = "3.12"; mv = "0.0487"; rank = 2 },
    @{ idx = 4; code = "010157"; name = "汇安中证500指数增强A";               size = "0.64";  pos = "86.25"; pct = "1.34"; mv = "0.0086"; rank = 7 },
    @{ idx = 5; code = "010158"; name = "汇安中证500指数增强C";               size = "0.46";  pos = "86.25"; pct = "1.34"; mv = "0.0062"; rank = 7 },
    @{ idx = 6; code = "016347"; name = "招商中证煤炭等权指数（LOF）E";        size = "0.20";  pos = "93.84"; pct = "3.12"; mv = "0.0062"; rank = 2 },
    @{ idx = 7; code = "009263"; name = "华宝红利精选混合A";                  size = "0.27";  pos = "91.60"; pct = "2.12"; mv = "0.0057"; rank = 1 },
    @{ idx = 8; code = "010841"; name = "华宝红利精选混合C";                  size = "0.18";  pos = "91.60"; pct = "2.12"; mv = "0.0038"; rank = 1 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.idx
    Set-TextCell $newSheet $r 2 $row.code
    Set-TextCell $newSheet $r 3 $row.name
    Set-TextCell $newSheet $r 4 $row.size
    Set-TextCell $newSheet $r 5 $row.pos
    Set-TextCell $newSheet $r 6 $row.pct
    Set-TextCell $newSheet $r 7 $row.mv
    $newSheet.Cells.Item($r, 8).Value = $row.rank
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 right
#    after the header row, pushing all existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summaryRows = @(
    @{ idx = 0; label = "2022-Q4"; count = 9;  mv = 0.88 },
    @{ idx = 1; label = "2022-Q3"; count = 19; mv = 1.02 },
    @{ idx = 2; label = "2022-Q2"; count = 4;  mv = 0.89 },
    @{ idx = 3; label = "2022-Q1"; count = 13; mv = 0.59 },
    @{ idx = 4; label = "2021-Q4"; count = 1;  mv = 0.07000000000000001 },
    @{ idx = 5; label = "2021-Q3"; count = 11; mv = 2.88 },
    @{ idx = 6; label = "2021-Q1"; count = 5;  mv = 0.15 },
    @{ idx = 7; label = "2020-Q4"; count = 2;  mv = 0.07000000000000001 }
)

$summary.Range("A3").Copy($summary.Range("A2"))

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row.idx
    $summary.Cells.Item($r, 2).Value = $row.label
    $summary.Cells.Item($r, 3).Value = $row.count
    $summary.Cells.Item($r, 4).Value = $row.mv
    $r = $r + 1
}
